# Updates TabQuery/StatQuery/TsvExcel/WebExcel cell contents for phs002430 Sequencing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tsvFileName = @'
TC03_CDS_phs002430_ExperimentalStrategy_Sequencing_TSVData.xlsx
'@

$webFileName = @'
TC03_CDS_phs002430_ExperimentalStrategy_Sequencing_WebData.xlsx
'@

$statQuery = @'
SELECT
    COUNT(DISTINCT s.study_name) AS "Studies",
    COUNT(DISTINCT sp.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    COUNT(DISTINCT f.file_id) AS "Files"
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN 
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN 
    df_file f ON f."sample.sample_id" = smp.sample_id
JOIN 
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN 
    df_genomic_info gi ON gi."file.file_id" = f.file_id
WHERE 
     s.phs_accession = 'phs002430' AND f.experimental_strategy_and_data_subtypes = 'Sequencing';
'@

$participantsQuery = @'
WITH Distinct_Samples AS (
    SELECT DISTINCT
        sp.participant_id,
        sp.study_participant_id,
        s.study_name,
        s.phs_accession,
        sp.gender,
        smp.sample_id
    FROM 
        df_participant sp
    JOIN 
        df_study s ON sp."study.phs_accession" = s.phs_accession
    JOIN 
        df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_program p ON p.program_acronym = s."program.program_acronym"
    JOIN
        df_file f1 ON f1."sample.sample_id" = smp.sample_id
    JOIN
        df_genomic_info gi ON gi."file.file_id" = f1.file_id
    WHERE 
        s.phs_accession = 'phs002430' AND f1.experimental_strategy_and_data_subtypes = 'Sequencing'
),
Sample_Limit AS (
    SELECT
        participant_id,
        study_participant_id,
        study_name,
        phs_accession,
        gender,
        sample_id,
        ROW_NUMBER() OVER (PARTITION BY study_participant_id ORDER BY sample_id) as rn
    FROM 
        Distinct_Samples
)
SELECT
    participant_id AS "Participant ID", 
    study_name AS "Study Name",
    phs_accession AS Accession,
    gender AS Gender,
    GROUP_CONCAT(
        CASE 
            WHEN rn <= 5 THEN sample_id 
        END, ', '
    ) ||
    CASE 
        WHEN MAX(rn) > 5 THEN ', ...' 
        ELSE '' 
    END AS Samples
FROM 
    Sample_Limit
GROUP BY
    participant_id, 
    study_name,
    phs_accession,
    gender
LIMIT 100;
'@

$samplesQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
     s.phs_accession = 'phs002430' AND f1.experimental_strategy_and_data_subtypes = 'Sequencing'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$filesQuery = @'
SELECT
    f1.file_name AS "File Name",
    s.study_name AS "Study Name",
    s.phs_accession AS "Accession",
    sp.participant_id AS "Participant Id",
    COALESCE((
        SELECT
            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||
            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END
        FROM (
            SELECT
                smp.sample_id,
                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn
            FROM df_sample smp
            WHERE smp."participant.study_participant_id" = sp.study_participant_id
        ) smp
    ), '') AS "Sample Id",
    f1.file_type AS "File Type",
    gi.library_strategy AS "Library Strategy"
FROM 
    df_study s
INNER JOIN 
    df_participant sp ON sp."study.phs_accession" = s.phs_accession
INNER JOIN  
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
INNER JOIN 
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
INNER JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
INNER JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
INNER JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
WHERE 
    s.phs_accession = 'phs002430' AND f1.experimental_strategy_and_data_subtypes = 'Sequencing'
GROUP BY
    f1.file_name,
    s.study_name,
    s.phs_accession,
    sp.participant_id,
    f1.file_type,
    gi.library_strategy
ORDER BY 
    f1.file_name ASC
LIMIT 100;
'@

$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $tsvFileName
$ws.Range("E2").Value = $webFileName
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

$ws.Rows.Item(2).RowHeight = 409.5

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C9").Select()
